# Update cryptos list data (simulating refreshed GitHub Actions scrape)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.386.32', '  -0.32%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.847.75', '  -0.13%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9988', '  -0.06%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '240.38', '  -0.59%  '),
    @(6, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.6306', '  +0.40%  '),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.000', '  -0.02%  '),
    @(8, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07562', '  +0.47%  '),
    @(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2959', '  -0.50%  '),
    @(10, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '24.44', '  +0.21%  '),
    @(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07702', '  -0.42%  '),
    @(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.867.02', '  +0.55%  '),
    @(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.999', '  -0.01%  '),
    @(14, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6839', '  -0.98%  '),
    @(15, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001004', '  +2.93%  '),
    @(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '82.96', '  -0.63%  '),
    @(17, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '2.122.49', '  +0.35%  '),
    @(18, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.154', '  -1.21%  '),
    @(19, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.410.77', '  -0.41%  '),
    @(20, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '227.86', '  -2.23%  '),
    @(21, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '12.47', '  -0.11%  '),
    @(22, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9997', '  -0.01%  '),
    @(23, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.564', '  -0.81%  '),
    @(24, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '1.001', '  -0.01%  '),
    @(25, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '157.54', '  +1.97%  '),
    @(26, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1397', '  +0.73%  '),
    @(27, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.382', '  -0.84%  '),
    @(28, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.66', '  -0.20%  '),
    @(29, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.464', '  -0.71%  '),
    @(30, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.265', '  +1.00%  '),
    @(31, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.05687', '  -3.92%  '),
    @(32, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.123', '  +0.42%  '),
    @(33, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.025', '  -0.05%  '),
    @(34, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.853', '  -1.23%  '),
    @(35, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.155', '  -1.15%  '),
    @(36, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.7133', '  -0.70%  '),
    @(37, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.589', '  -0.05%  '),
    @(38, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.262.78', '  +1.93%  '),
    @(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01815', '  +1.33%  '),
    @(40, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.779', '  -0.53%  '),
    @(41, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.9060', '  -0.09%  '),
    @(42, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.169', '  +0.83%  '),
    @(43, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9999', '  +0.00%  '),
    @(44, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '101.51', '  +0.18%  '),
    @(45, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '66.34', '  -1.26%  '),
    @(46, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.076', '  -4.14%  '),
    @(47, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.4045', '  +0.16%  '),
    @(48, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.101', '  -0.76%  '),
    @(49, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.684', '  -0.72%  '),
    @(50, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1121', '  +0.27%  '),
    @(51, 'Cronos', 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro', '0.05734', '  -0.36%  ')
)

foreach ($entry in $data) {
    $r = $entry[0]
    $ws.Cells.Item($r, 2).Value = $entry[1]
    $ws.Cells.Item($r, 3).Value = $entry[2]
    $ws.Cells.Item($r, 4).Value = $entry[3]
    $ws.Cells.Item($r, 5).Value = $entry[4]
}
